$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the "-Dean's List" paragraph to
#    the end of the objective paragraph. Remove it from its current
#    location first (it's a hidden bookmark, so it must be looked up
#    by name rather than enumerated).
# -----------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# -----------------------------------------------------------------
# 2) Rewrite the objective paragraph ("To obtain a position working
#    with embedded systems ...") as a sequence of runs matching the
#    target revision, then re-add the "_GoBack" bookmark at its end.
# -----------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "To obtain a position working with embedded systems that leverages my skills, experience, and accomplishments as a developer."
$find.Forward = $true
$find.Wrap = 0
$found = $find.Execute()

if (-not $found) {
    throw "Objective paragraph text not found"
}

$target = $d.Content
$target.Start = $find.Parent.Start
$target.End = $find.Parent.End

$target.Text = ""

$runsXml = ""
$runsXml += '<w:r><w:t xml:space="preserve">To obtain a position working with </w:t></w:r>'
$runsXml += '<w:r><w:t>an Agile team</w:t></w:r>'
$runsXml += '<w:r><w:t xml:space="preserve"> that lever</w:t></w:r>'
$runsXml += '<w:r><w:t xml:space="preserve">ages my skills, experience, and </w:t></w:r>'
$runsXml += '<w:r><w:t>accomplishments as a developer</w:t></w:r>'
$runsXml += '<w:r><w:t xml:space="preserve"> while providing me with new skills to succeed</w:t></w:r>'
$runsXml += '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
$runsXml += '<w:r><w:t>with</w:t></w:r>'
$runsXml += '<w:r><w:t xml:space="preserve">in my career. </w:t></w:r>'
$runsXml += '<w:bookmarkStart w:id="0" w:name="_GoBack"/>'
$runsXml += '<w:bookmarkEnd w:id="0"/>'

$pkgXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint = $d.Range($target.Start, $target.Start)
$insertionPoint.InsertXML($pkgXml)

Write-Output "objective paragraph updated"
